# Update to 2021_Q1 data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Project Achilles: common genes" row (row 2) values ---
$ws.Range("B2").Value = 18119
$ws.Range("C2").Value = 2085
$ws.Range("E2").Value = 15673
$ws.Range("F2").Value = 2010
$ws.Range("I2").Value = 1.6199226075875801
$ws.Range("J2").Value = 1.24048204990668
$ws.Range("K2").Value = 2.09234944004079
$ws.Range("L2").Value = 0.00033745901166429298

# --- Refresh formatting for the data rows (rows 2-7), bumping the font
#     size back to the workbook default (11pt) while preserving bold
#     weight per column, which is what produces the new style entries. ---
$ws.Range("B2:G7").Font.Size = 11
$ws.Range("H2:H7").Font.Size = 11
$ws.Range("I2:I7").Font.Size = 11
$ws.Range("J2:K7").Font.Size = 11
$ws.Range("L2:L7").Font.Size = 11
